$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, G (column F is unchanged) for rows 2-7.
$data = @{
    2 = @(0.3464964993005633, 9.226618575922256, 16.98373111632243, 6.48142807727062, 33.03827426881587)
    3 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 6.48142807727062, 28.30127388105354)
    4 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    5 = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.811642989160245)
    6 = @(0.02258322285507441, 0.3375848360084654, 0.1529057820181812, 246.9852506941017, 247.4983245349834)
    7 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 6.741336633845642)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
